# Automatische test-sync: 2025-08-19 20:24:50
# Append a new log row to the "Logs" sheet and update the "Dashboard" count.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$newRow = 16

$ws.Range("A$newRow").Value = "Interne taak"
$ws.Range("B$newRow").Value = "kwaliteit@testbedrijf123.nl"
$ws.Range("D$newRow").Value = "Intern verzoek / Actie voor medewerker"
$ws.Range("F$newRow").Value = "2025-08-19 20:24:05"
$ws.Range("G$newRow").Value = "Nee"
$ws.Range("H$newRow").Value = "Ja"
$ws.Range("I$newRow").Value = "Nee"
$ws.Range("J$newRow").Value = "Nee"

# Extend the conditional-formatting "applies to" ranges (D/G/H/I/J) so the
# new row 16 is covered too, same as the original row 2:15 -> 2:16 bump.
$ws.Range("D2:D$newRow").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D$newRow"))
$ws.Range("G2:G$newRow").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G$newRow"))
$ws.Range("H2:H$newRow").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H$newRow"))
$ws.Range("I2:I$newRow").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I$newRow"))
$ws.Range("J2:J$newRow").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J$newRow"))

# Bump the Dashboard summary count to reflect the new log row.
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 15

Write-Output "applied"
